# Add the new "Program Area" worksheet after the existing "Client Time" sheet,
# populate it with its lookup table, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Program Area"

# Header row
$ws.Range("A1").Value = "ProgramArea"
$ws.Range("B1").Value = "Value"

# Program area rows
$ws.Range("A2").Value = "PrEP_NEW"
$ws.Range("A3").Value = "PrEP_CURR"
$ws.Range("A4").Value = "HTS_SELF"
$ws.Range("A5").Value = "HTS_TST"
$ws.Range("A6").Value = "TX_NEW"
$ws.Range("A7").Value = "TX_CURR"
$ws.Range("A8").Value = "PMTCT_ART"
$ws.Range("A9").Value = "TX_PVLS"

# Stray formatted (but empty) cell next to the second row
$ws.Range("C2").WrapText = $true

# Column A sized to fit the longest label
$ws.Columns.Item(1).ColumnWidth = 14.2

# Turn the range into a table, matching the style used elsewhere in the workbook
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:B9"), $null, 1)
$tbl.Name = "Table2"
$tbl.TableStyle = "TableStyleMedium1"

# Leave the selection where the author left it and make this the active tab
$ws.Range("F13").Select() | Out-Null
$ws.Activate()
